$d = $word.ActiveDocument

# Locate the paragraph that carries the astromap credit / link text.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*Os mapas de*GaNight*") {
        $target = $p
        break
    }
}

$r = $target.Range
$start = $r.Start
$end = $r.End - 1   # exclude the paragraph mark

# Remove all of the existing (heavily run-split) paragraph content.
$old = $d.Range($start, $end)
$old.Delete()

# New, single-run text with the year bumped from 2018 to 2022.
$newText = "Os mapas de estrelas deste documento foron preparados por Jenik Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."

$insPoint = $d.Range($start, $start)
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
       '<pkg:xmlData>' +
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:body><w:p><w:r/><w:r><w:t>' + $newText + '</w:t></w:r></w:p></w:body>' +
       '</w:document>' +
       '</pkg:xmlData></pkg:part></pkg:package>'
$insPoint.InsertXML($xml)
